# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to match the freshly scraped numbers (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Per-sheet (F column, and occasionally G column) new values, keyed by row.
# 展览 (sheet1) and 全部类型 (sheet4) are independent copies of the same
# scrape, so a couple of rows (12, 37, 45) pick up slightly different counts
# between the two sheets.
$sheetUpdates = @{
    "展览" = @{
        3  = @{ F = 1460 }
        9  = @{ F = 199 }
        12 = @{ F = 4838 }
        14 = @{ F = 7125 }
        21 = @{ F = 4211 }
        22 = @{ F = 1524 }
        25 = @{ F = 2785 }
        29 = @{ F = 410 }
        33 = @{ F = 61 }
        35 = @{ F = 1096 }
        36 = @{ F = 79 }
        37 = @{ F = 1120 }
        38 = @{ F = 94 }
        41 = @{ F = 504 }
        43 = @{ F = 32 }
        45 = @{ F = 1685; G = 65 }
        46 = @{ F = 667 }
    }
    "全部类型" = @{
        3  = @{ F = 1460 }
        9  = @{ F = 199 }
        12 = @{ F = 4839 }
        14 = @{ F = 7125 }
        21 = @{ F = 4211 }
        22 = @{ F = 1524 }
        25 = @{ F = 2785 }
        29 = @{ F = 410 }
        33 = @{ F = 61 }
        35 = @{ F = 1096 }
        36 = @{ F = 79 }
        37 = @{ F = 1121 }
        38 = @{ F = 94 }
        41 = @{ F = 504 }
        43 = @{ F = 32 }
        45 = @{ F = 1686; G = 65 }
        46 = @{ F = 667 }
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cols = $rows[$rowNum]
        foreach ($colLetter in $cols.Keys) {
            $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
        }
    }
}
